$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 183-184, pushing the existing data
# (old rows 183-218) down to rows 185-220.
$ws.Range("A183:A184").EntireRow.Insert()

# New row 183 data
$ws.Cells.Item(183, 1).Value = 10
$ws.Cells.Item(183, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(183, 3).Value = 'La Araucanía'
$ws.Cells.Item(183, 4).Value = 44476
$ws.Cells.Item(183, 5).Value = 9
$ws.Cells.Item(183, 6).Value = 'Fruta'
$ws.Cells.Item(183, 7).Value = 100108
$ws.Cells.Item(183, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(183, 9).Value = 100108002
$ws.Cells.Item(183, 10).Value = 'Mango'
$ws.Cells.Item(183, 11).Value = 'Sin especificar'
$ws.Cells.Item(183, 12).Value = 'Primera'
$ws.Cells.Item(183, 13).Value = 1200
$ws.Cells.Item(183, 14).Value = 8000
$ws.Cells.Item(183, 15).Value = 9000
$ws.Cells.Item(183, 16).Value = 8417
$ws.Cells.Item(183, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(183, 18).Value = 'Perú'
$ws.Cells.Item(183, 19).Value = 2104
$ws.Cells.Item(183, 20).Value = 4

# New row 184 data
$ws.Cells.Item(184, 1).Value = 10
$ws.Cells.Item(184, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(184, 3).Value = 'La Araucanía'
$ws.Cells.Item(184, 4).Value = 44476
$ws.Cells.Item(184, 5).Value = 9
$ws.Cells.Item(184, 6).Value = 'Fruta'
$ws.Cells.Item(184, 7).Value = 100108
$ws.Cells.Item(184, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(184, 9).Value = 100108002
$ws.Cells.Item(184, 10).Value = 'Mango'
$ws.Cells.Item(184, 11).Value = 'Sin especificar'
$ws.Cells.Item(184, 12).Value = 'Segunda'
$ws.Cells.Item(184, 13).Value = 200
$ws.Cells.Item(184, 14).Value = 7000
$ws.Cells.Item(184, 15).Value = 7000
$ws.Cells.Item(184, 16).Value = 7000
$ws.Cells.Item(184, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(184, 18).Value = 'Perú'
$ws.Cells.Item(184, 19).Value = 1750
$ws.Cells.Item(184, 20).Value = 4

# Make sure date cells use the same numeric date style as the rest of column D
$ws.Range("D183:D184").NumberFormat = $ws.Range("D185").NumberFormat
